$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 18:44:45"
$ws.Range("A3").Value = "Total filas: 316"
$ws.Range("C64").Value = "15X38_ABASTO"
$ws.Range("C65").Value = "14_ABASTO"
$ws.Range("A102").Value = "09:22:34"
$ws.Range("C102").Value = "15X38_ABASTO"
$ws.Range("D102").Value = 113
$ws.Range("A103").Value = "10:56:15"
$ws.Range("C103").Value = "14_ABASTO"
$ws.Range("D103").Value = 19
$ws.Range("A133").Value = "10:36:50"
$ws.Range("C133").Value = "14_ABASTO"
$ws.Range("D133").Value = 116
$ws.Range("A134").Value = "11:33:52"
$ws.Range("C134").Value = "23_HERNANDEZ"
$ws.Range("D134").Value = 59
$ws.Range("A301").Value = "17:35:41"
$ws.Range("C301").Value = "27_EL RETIRO"
$ws.Range("D301").Value = 101
$ws.Range("A302").Value = "17:55:25"
$ws.Range("C302").Value = "17_ROMERO"
$ws.Range("D302").Value = 81
$ws.Range("A304").Value = "18:44:45"
$ws.Range("B304").Value = "19:21"
$ws.Range("D304").Value = 37
$ws.Range("A305").Value = "17:55:25"
$ws.Range("B305").Value = "19:22"
$ws.Range("C305").Value = "23_HERNANDEZ"
$ws.Range("D305").Value = 87
$ws.Range("B306").Value = "19:23"
$ws.Range("C306").Value = "16_SANTA ANA"
$ws.Range("D306").Value = 53
$ws.Range("A307").Value = "18:30:48"
$ws.Range("B307").Value = "19:25"
$ws.Range("D307").Value = 55
$ws.Range("A308").Value = "18:11:09"
$ws.Range("C308").Value = "23_HERNANDEZ"
$ws.Range("D308").Value = 76
$ws.Range("A309").Value = "17:35:41"
$ws.Range("B309").Value = "19:27"
$ws.Range("D309").Value = 112
$ws.Range("A310").Value = "17:55:25"
$ws.Range("B310").Value = "19:28"
$ws.Range("C310").Value = "215C_EL PATO"
$ws.Range("D310").Value = 93
$ws.Range("A311").Value = "17:47:45"
$ws.Range("B311").Value = "19:35"
$ws.Range("D311").Value = 108
$ws.Range("B312").Value = "19:36"
$ws.Range("C312").Value = "11_ETCHEVERRY"
$ws.Range("D312").Value = 101
$ws.Range("A313").Value = "17:55:25"
$ws.Range("B313").Value = "19:39"
$ws.Range("D313").Value = 104
$ws.Range("A314").Value = "17:47:45"
$ws.Range("B314").Value = "19:42"
$ws.Range("C314").Value = "15X38_ABASTO"
$ws.Range("D314").Value = 115
$ws.Range("B315").Value = "19:52"
$ws.Range("C315").Value = "81_EL PELIGRO"
$ws.Range("D315").Value = 117
$ws.Range("A316").Value = "17:55:25"
$ws.Range("B316").Value = "19:53"
$ws.Range("C316").Value = "225_GOMEZ"
$ws.Range("D316").Value = 118
$ws.Range("A317").Value = "18:11:09"
$ws.Range("B317").Value = "20:06"
$ws.Range("C317").Value = "215C_EL PATO"
$ws.Range("D317").Value = 115
$ws.Range("A318").Value = "18:44:45"
$ws.Range("B318").Value = "20:10"
$ws.Range("C318").Value = "23_HERNANDEZ"
$ws.Range("D318").Value = 86
$ws.Range("E318").Value = "LP1912"
$ws.Range("A319").Value = "18:44:45"
$ws.Range("B319").Value = "20:12"
$ws.Range("C319").Value = "14_ABASTO"
$ws.Range("D319").Value = 88
$ws.Range("E319").Value = "LP1912"
$ws.Range("A320").Value = "18:30:48"
$ws.Range("B320").Value = "20:21"
$ws.Range("C320").Value = "15_ABASTO"
$ws.Range("D320").Value = 111
$ws.Range("E320").Value = "LP1912"
$ws.Range("A321").Value = "18:44:45"
$ws.Range("B321").Value = "20:30"
$ws.Range("C321").Value = "10_OLMOS"
$ws.Range("D321").Value = 106
$ws.Range("E321").Value = "LP1912"

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 18:44:45"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 18:44:45"
